$wb = $excel.ActiveWorkbook

# Update the "Latest Handoff Datetime" for the last report row (6f40ab55...) in
# the zh-cn sheet to reflect the new handoff generated for this report.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-02-18 07:29:20"

# Same update for the de-de sheet.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-02-18 07:29:31"
